# Update the date placeholder on slide 1
$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(3).TextFrame.TextRange.Text = "2024-07-24"

# Merge the two paragraphs on slide 3 into a single paragraph/run
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "El número de datos es de 19452."
